$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '69.401.91'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -0.42%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.708.63'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +8.89%  '
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '610.39'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +4.02%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '174.99'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -4.14%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '3.705.16'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +9.25%  '
$ws.Range("E8").Value = '  +0.05%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.536'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.99%  '
$ws.Range("E10").Value = '  +4.61%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '6.34'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -3.42%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.495'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +1.73%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '40.48'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +5.93%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.0000252'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +0.76%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '4.329.24'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +8.57%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.707.24'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +9.61%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '69.445.48'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("E18").Value = '  +0.28%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '7.53'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +1.61%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '510.84'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +1.86%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '16.60'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -2.68%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '9.35'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +13.06%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.723'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -0.49%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '87.38'
$cell.Style = "Normal"
$ws.Range("E24").Value = '  +1.39%  '
$ws.Range("E25").Value = '  +2.78%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '13.33'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +0.23%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '10.92'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -0.32%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.0000124'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  +23.07%  '
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '2.49'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +0.60%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '2.83'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +5.07%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '7.81'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -4.91%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '31.02'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +4.27%  '
$ws.Range("E34").Value = '  -1.14%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.997'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -0.14%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.14'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +2.10%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '1.03'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +3.68%  '
$ws.Range("E38").Value = '  +2.19%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.17'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +2.90%  '
$ws.Range("E40").Value = '  +2.73%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '51.17'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +2.15%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '8.75'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +1.17%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '43.91'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -10.10%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '3.068.95'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +4.29%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '415.02'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.37%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.69'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -4.45%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.0362'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.53%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '27.64'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.16%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("B50").Value = 'ThetaToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '2.49'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +2.22%  '
$ws.Range("B51").Value = 'Monero'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '134.63'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -1.01%  '
